# Applies the "show_san_input" related data update to the Build Room /
# EUC_Perth_Assets workbook: refreshed item counts on the two *_Items
# sheets, plus newly logged SAN-add timestamp rows on the two
# *_Timestamps sheets.

$wb = $excel.ActiveWorkbook

# Helper: write a value that must be stored as TEXT even when it looks
# like a pure number (e.g. "12345"). Excel's normal Value setter treats
# digit-only strings as numbers, so the SAN number is entered as a
# formula yielding that literal string and then "baked" back into a
# plain value in place (copy + paste-values) - this keeps the cell's
# type as text without touching/creating any cell formatting/style.
function Set-TextValue($cell, $text) {
    $cell.Formula = "=""" + $text + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# 4.2_Items - refreshed LastCount / NewCount for three rows
# ---------------------------------------------------------------------
$wsItems = $wb.Worksheets.Item("4.2_Items")
$wsItems.Cells.Item(3, 2).Value = 1    # Dock Thunderbolt G4 - LastCount
$wsItems.Cells.Item(3, 3).Value = 21   # Dock Thunderbolt G4 - NewCount
$wsItems.Cells.Item(4, 2).Value = 0    # Laptop 840 G10 - LastCount
$wsItems.Cells.Item(4, 3).Value = 3    # Laptop 840 G10 - NewCount
$wsItems.Cells.Item(5, 2).Value = 1    # Laptop 840 G9 - LastCount
$wsItems.Cells.Item(5, 3).Value = 4    # Laptop 840 G9 - NewCount

# ---------------------------------------------------------------------
# 4.2_Timestamps - append newly logged "Add 1" SAN events (rows 85-89)
# ---------------------------------------------------------------------
$wsTs = $wb.Worksheets.Item("4.2_Timestamps")

$wsTs.Cells.Item(85, 1).Value = "2023-12-11 20:28:29"
$wsTs.Cells.Item(85, 2).Value = "Laptop 840 G9"
$wsTs.Cells.Item(85, 3).Value = "Add 1"
Set-TextValue $wsTs.Cells.Item(85, 4) "12345"

$wsTs.Cells.Item(86, 1).Value = "2023-12-11 20:28:32"
$wsTs.Cells.Item(86, 2).Value = "Laptop 840 G9"
$wsTs.Cells.Item(86, 3).Value = "Add 1"
Set-TextValue $wsTs.Cells.Item(86, 4) "55555"

$wsTs.Cells.Item(87, 1).Value = "2023-12-11 20:32:27"
$wsTs.Cells.Item(87, 2).Value = "Laptop 840 G10"
$wsTs.Cells.Item(87, 3).Value = "Add 1"
Set-TextValue $wsTs.Cells.Item(87, 4) "33333"

$wsTs.Cells.Item(88, 1).Value = "2023-12-11 20:32:33"
$wsTs.Cells.Item(88, 2).Value = "Laptop 840 G10"
$wsTs.Cells.Item(88, 3).Value = "Add 1"
Set-TextValue $wsTs.Cells.Item(88, 4) "33333"

$wsTs.Cells.Item(89, 1).Value = "2023-12-11 20:32:45"
$wsTs.Cells.Item(89, 2).Value = "Laptop 840 G10"
$wsTs.Cells.Item(89, 3).Value = "Add 1"
Set-TextValue $wsTs.Cells.Item(89, 4) "44444"

# ---------------------------------------------------------------------
# BR_Items - refreshed LastCount / NewCount for two rows
# ---------------------------------------------------------------------
$wsBrItems = $wb.Worksheets.Item("BR_Items")
$wsBrItems.Cells.Item(4, 2).Value = 0  # Laptop 840 G10 - LastCount
$wsBrItems.Cells.Item(4, 3).Value = 2  # Laptop 840 G10 - NewCount
$wsBrItems.Cells.Item(5, 2).Value = 0  # Laptop 840 G9 - LastCount
$wsBrItems.Cells.Item(5, 3).Value = 2  # Laptop 840 G9 - NewCount

# ---------------------------------------------------------------------
# BR_Timestamps - append the newly logged "Add 1" SAN event (row 17)
# ---------------------------------------------------------------------
$wsBrTs = $wb.Worksheets.Item("BR_Timestamps")
$wsBrTs.Cells.Item(17, 1).Value = "2023-12-11 20:28:50"
$wsBrTs.Cells.Item(17, 2).Value = "Laptop 840 G10"
$wsBrTs.Cells.Item(17, 3).Value = "Add 1"
Set-TextValue $wsBrTs.Cells.Item(17, 4) "22222"
